# basic 2-Dim Arrray example 6.3.4EKMEDJE
# Adds a new "Z6 - LOOPs" section (header + sum row + 31 data rows) below the
# existing "Z5 - REVISION" section, updates the TOTAL formulas in row 2 to
# include the new section, and bolds/shades the "Z5 - REVISION" header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-style the existing "Z5 - REVISION" header (row 80): bold + shaded fill
#    (previously just centered with no fill/bold).
# ---------------------------------------------------------------------------
$c80 = $ws.Range("C80")
$d80 = $ws.Range("D80")

$c69 = $ws.Range("C69")   # already has the fillId2 + centered look we want to reuse
$c69.Copy()
$c80.PasteSpecial(-4122)  # xlPasteFormats
$c69.Copy()
$d80.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$c80.Font.Bold = $true
$d80.Font.Bold = $true

# ---------------------------------------------------------------------------
# 2) New section header row 82: "Z6 - LOOPs", merged C82:D82, same look as the
#    other "SUM row" shaded+centered style (fillId2 + center), reusing C69's
#    formatting again (not bold).
# ---------------------------------------------------------------------------
$c82 = $ws.Range("C82")
$d82 = $ws.Range("D82")

$c69.Copy()
$c82.PasteSpecial(-4122)
$c69.Copy()
$d82.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$c82.Value = "Z6 - LOOPs"
$ws.Range("C82:D82").Merge()

# ---------------------------------------------------------------------------
# 3) New SUM row 83: plain shaded (fillId2, no center) like row 5 (Zadanie 1).
# ---------------------------------------------------------------------------
$c83 = $ws.Range("C83")
$d83 = $ws.Range("D83")

$c5 = $ws.Range("C5")
$c5.Copy()
$c83.PasteSpecial(-4122)
$c5.Copy()
$d83.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$c83.Formula = "=SUM(C84:C114)"
$d83.Formula = "=SUM(D84:D114)"

# ---------------------------------------------------------------------------
# 4) 31 raw data rows (84-114), no special style (same as other data blocks).
# ---------------------------------------------------------------------------
$data = @(
  @(84, 4, 7),
  @(85, 6, 31),
  @(86, 9, 32),
  @(87, 12, 43),
  @(88, 9, 25),
  @(89, 9, 38),
  @(90, 3, 4),
  @(91, 7, 7),
  @(92, 14, 41),
  @(93, 2, 20),
  @(94, 2, 17),
  @(95, 50, 58),
  @(96, 44, 39),
  @(97, 7, 38),
  @(98, 4, 52),
  @(99, 41, 22),
  @(100, 44, 39),
  @(101, 41, 11),
  @(102, 60, 2),
  @(103, 37, 21),
  @(104, 3, 1),
  @(105, 29, 9),
  @(106, 24, 19),
  @(107, 31, 26),
  @(108, 40, 49),
  @(109, 58, 35),
  @(110, 5, 6),
  @(111, 32, 35),
  @(112, 21, 22),
  @(113, 42, 6),
  @(114, 11, 38)
)

foreach ($row in $data) {
    $r = $row[0]
    $cval = $row[1]
    $dval = $row[2]
    $ws.Cells.Item($r, 3).Value = $cval
    $ws.Cells.Item($r, 4).Value = $dval
}

# ---------------------------------------------------------------------------
# 5) Update the TOTAL formulas in row 2 so they roll up the new Z6 section.
# ---------------------------------------------------------------------------
$ws.Range("C2").Formula = "=SUM(C5+C13+C24+C69+C81+C83)+(D2/60)"
$ws.Range("D2").Formula = "=SUM(D5+D24+D69+D81++C81+C83)"
$ws.Range("A2").Formula = "=C2/80"
$ws.Range("B2").Formula = "=C2/45"

# ---------------------------------------------------------------------------
# 6) Scroll / selection bookkeeping to mirror the author's view state.
# ---------------------------------------------------------------------------
$ws.Range("D84").Select()

$wb.Save()
